# -----------------------------------------------------------------------------
# Updates cryptos list (price / 1h-volume columns, plus a few re-ranked rows)
# as published by the "Updated cryptos list ... with GitHub Actions" commit.
# -----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.791.53"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "3.808.59"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'597.60"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "'167.78"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "'0.160"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").Value = "'6.29"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "'36.02"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "4.447.35"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").Value = "3.790.20"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "'18.55"
$ws.Range("E16").Value = "  +5.34%  "
$ws.Range("D17").Value = "67.797.28"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'461.83"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "'9.95"
$ws.Range("E21").Value = "  -5.08%  "
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("D24").Value = "'83.64"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "'12.10"
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("D27").Value = "'10.02"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "3.957.68"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  +4.96%  "
$ws.Range("D32").Value = "'7.30"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").Value = "'29.69"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "3.748.47"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "'0.100"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").Value = "'3.43"
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D44").Value = "'48.15"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "'8.33"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'148.54"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").Value = "'396.33"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'26.98"
$ws.Range("E50").Value = "  +7.91%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.84"
$ws.Range("E51").Value = "  -2.13%  "
